# Update metadata (Metadata sheet) and Elements sheet per commit
# "Update description et fix compilation"

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-ballot -> 2.0.0
$ws1.Range("B3").Value = "2.0.0"

# Date: refreshed publication timestamp
$ws1.Range("B8").Value = "2025-10-16T14:47:34+00:00"

# Description: richer, itemised French description
$description = @"

- **Séjour** : commentaire relatif au séjour.
- **Événement** : commentaires sur le déroulé de l’évènement.
- **Évaluation** : commentaire libre sur le contenu ou le résultat de l’évaluation.
- **Champ évalué** : commentaire spécifique à un item ou sous-item évalué.
"@
$description = $description.TrimEnd("`r", "`n")
$ws1.Range("B12").Value = $description

# Remove the extra Context row (element:QuestionnaireResponse), row 22
$ws1.Rows.Item(22).Delete()

# --- Elements sheet ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Elements")

# The "Extension" row's Definition (column M, row 2) mirrors the same
# text as the Metadata "Description" - keep both in sync.
$ws2.Range("M2").Value = $description

# Extension.value[x] Short column (L, row 6): replace generic text with
# a concrete French example comment.
$ws2.Range("L6").Value = "Exemple de commentaire : Cet évènement a débuté plus tard l’usager était sous la douche à l’heure du début du rendez-vous."
